# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns.
# Numeric-looking Price strings are prefixed with a leading apostrophe so
# Excel keeps them as text (matching the source data) instead of coercing
# them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.426.49"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "3.621.76"
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'603.31"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "'196.24"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("E7").Value = "  -1.00%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("D10").Value = "'0.645"
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("D11").Value = "'53.21"
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").Value = "'0.0000302"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").Value = "'9.57"
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").Value = "4.200.60"
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("D15").Value = "'600.67"
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("D16").Value = "'12.94"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").Value = "70.571.81"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "3.607.32"
$ws.Range("E18").Value = "  +1.80%  "
$ws.Range("D19").Value = "'19.04"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("D21").Value = "'0.996"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "'18.49"
$ws.Range("E22").Value = "  +2.12%  "
$ws.Range("D23").Value = "'5.21"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("D24").Value = "'102.46"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  -4.30%  "
$ws.Range("D27").Value = "'10.61"
$ws.Range("E27").Value = "  -2.73%  "
$ws.Range("D28").Value = "'9.67"
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("D29").Value = "'33.77"
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("D30").Value = "'4.66"
$ws.Range("E30").Value = "  +8.17%  "
$ws.Range("D31").Value = "'7.33"
$ws.Range("E31").Value = "  +2.96%  "
$ws.Range("D32").Value = "'12.27"
$ws.Range("E32").Value = "  -2.82%  "
$ws.Range("E33").Value = "  +2.29%  "
$ws.Range("D34").Value = "'63.45"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("D35").Value = "0.0₃0883"
$ws.Range("E35").Value = "  +2.94%  "
$ws.Range("D36").Value = "3.928.16"
$ws.Range("E36").Value = "  +5.48%  "
$ws.Range("D37").Value = "'533.08"
$ws.Range("E37").Value = "  +8.38%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("D39").Value = "'3.06"
$ws.Range("E39").Value = "  +1.05%  "
$ws.Range("D40").Value = "'36.81"
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("D42").Value = "'3.53"
$ws.Range("E42").Value = "  -3.44%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").Value = "'0.0462"
$ws.Range("E44").Value = "  +1.29%  "
$ws.Range("D45").Value = "'3.48"
$ws.Range("E45").Value = "  +5.37%  "
$ws.Range("D46").Value = "'2.87"
$ws.Range("E46").Value = "  +1.32%  "
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").Value = "'8.57"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("E51").Value = "  +1.62%  "
